$wb = $excel.ActiveWorkbook

# Map of sheet name -> { row -> new value for column F }
$updates = @{
    "展览" = @{
        2  = 3158
        4  = 1110
        8  = 46
        9  = 1142
        10 = 15996
        11 = 258
        14 = 6255
        16 = 117
        26 = 214
        31 = 11187
        35 = 187
    }
    "全部类型" = @{
        2  = 3158
        4  = 1110
        8  = 46
        9  = 1142
        10 = 15996
        11 = 258
        14 = 6255
        16 = 117
        26 = 214
        32 = 11187
        36 = 187
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Range("F$row").Value = $rows[$row]
    }
}
